$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1 (09:00 -> 09:15)
$ws.Range("F1").Value = "Last status check on: 20.01.2022 09:15"

# Row 8 (Benzina Albert Modrice) was refreshed with a new price check:
#  - current price (B8) and previous price (C8) swapped values
#  - delta (D8) is now stored as the text "+0.3" instead of the number -0.3
#  - last-checked date (E8) is now stored as plain text instead of a date serial
$ws.Range("B8").Value = 36.2
$ws.Range("C8").Value = 35.9

$ws.Range("D8").Value = "'+0.3"
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = "2022-01-20 09:15:18"
$ws.Range("E8").Style = "Normal"
